# ---------------------------------------------------------------------------
# Edit: MEGATAB_EMPREEND_JUN2025vlight.xlsx
#
# Splits the combined "COORDENADA(DEC)" column (E) - which stored strings
# like "-38.455678,-3.891234,0" (lon,lat,0) - into two new trailing columns:
#   N = LATITUDE
#   O = LONGITUDE
# The old column E header becomes "VER NO MAPA" and its data cells (E2:E23)
# are cleared. The AutoFilter / _FilterDatabase range collapses to the
# header row only (A1:O1) and the used range grows to A1:O23.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row: rename the coordinate column header, add new headers ---
$ws.Range("E1").Value = "VER NO MAPA"
$ws.Range("N1").Value = "LATITUDE"
$ws.Range("O1").Value = "LONGITUDE"

# --- 2. Remove the old combined coordinate values from column E ---
$ws.Range("E2:E23").ClearContents()

# --- 3. Stage the split latitude/longitude text values in helper columns
#        P:Q using TEXT() so they land as plain text (not re-parsed as
#        numbers), then copy/paste-special as values into N:O and clean up
#        the helper columns. ---
$ws.Range("P2").Formula = '=TEXT(-3.891234,"0.000000")'
$ws.Range("Q2").Formula = '=TEXT(-38.455678,"0.000000")'
$ws.Range("P3").Formula = '=TEXT(-3.789012,"0.000000")'
$ws.Range("Q3").Formula = '=TEXT(-38.512345,"0.000000")'
$ws.Range("P4").Formula = '=TEXT(-3.715678,"0.000000")'
$ws.Range("Q4").Formula = '=TEXT(-38.567890,"0.000000")'
$ws.Range("P5").Formula = '=TEXT(-3.812345,"0.000000")'
$ws.Range("Q5").Formula = '=TEXT(-38.543210,"0.000000")'
$ws.Range("P6").Formula = '=TEXT(-3.812345,"0.000000")'
$ws.Range("Q6").Formula = '=TEXT(-38.543210,"0.000000")'
$ws.Range("P7").Formula = '=TEXT(-3.732456,"0.000000")'
$ws.Range("Q7").Formula = '=TEXT(-38.489123,"0.000000")'
$ws.Range("P8").Formula = '=TEXT(-3.812345,"0.000000")'
$ws.Range("Q8").Formula = '=TEXT(-38.543210,"0.000000")'
$ws.Range("P9").Formula = '=TEXT(-3.732456,"0.000000")'
$ws.Range("Q9").Formula = '=TEXT(-38.489123,"0.000000")'
$ws.Range("P10").Formula = '=TEXT(-3.812345,"0.000000")'
$ws.Range("Q10").Formula = '=TEXT(-38.543210,"0.000000")'
$ws.Range("P11").Formula = '=TEXT(-3.824532,"0.000000")'
$ws.Range("Q11").Formula = '=TEXT(-38.579120,"0.000000")'
$ws.Range("P12").Formula = '=TEXT(-3.793874,"0.000000")'
$ws.Range("Q12").Formula = '=TEXT(-38.481262,"0.000000")'
$ws.Range("P13").Formula = '=TEXT(-3.816781,"0.000000")'
$ws.Range("Q13").Formula = '=TEXT(-38.551234,"0.000000")'
$ws.Range("P14").Formula = '=TEXT(-3.727890,"0.000000")'
$ws.Range("Q14").Formula = '=TEXT(-38.639012,"0.000000")'
$ws.Range("P15").Formula = '=TEXT(-3.718389,"0.000000")'
$ws.Range("Q15").Formula = '=TEXT(-38.482273,"0.000000")'
$ws.Range("P16").Formula = '=TEXT(-3.830000,"0.000000")'
$ws.Range("Q16").Formula = '=TEXT(-38.550000,"0.000000")'
$ws.Range("P17").Formula = '=TEXT(-3.714701,"0.000000")'
$ws.Range("Q17").Formula = '=TEXT(-38.581138,"0.000000")'
$ws.Range("P18").Formula = '=TEXT(-3.837602,"0.000000")'
$ws.Range("Q18").Formula = '=TEXT(-38.460851,"0.000000")'
$ws.Range("P19").Formula = '=TEXT(-3.812563,"0.000000")'
$ws.Range("Q19").Formula = '=TEXT(-38.537415,"0.000000")'
$ws.Range("P20").Formula = '=TEXT(-3.732028,"0.000000")'
$ws.Range("Q20").Formula = '=TEXT(-38.462216,"0.000000")'
$ws.Range("P21").Formula = '=TEXT(-3.892758,"0.000000")'
$ws.Range("Q21").Formula = '=TEXT(-38.455388,"0.000000")'
$ws.Range("P22").Formula = '=TEXT(-3.873219,"0.000000")'
$ws.Range("Q22").Formula = '=TEXT(-38.635111,"0.000000")'
$ws.Range("P23").Formula = '=TEXT(-3.727890,"0.000000")'
$ws.Range("Q23").Formula = '=TEXT(-38.639012,"0.000000")'

$stage = $ws.Range("P2:Q23")
$stage.Copy()
$ws.Range("N2").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$stage.Clear()

# --- 4. Reset the AutoFilter so it spans the new header-only range A1:O1 ---
$ws.AutoFilterMode = $false
$ws.Range("A1:O1").AutoFilter()

# --- 5. Point the _FilterDatabase defined name at the new range ---
$fd = $wb.Names.Item("Planilha1!_FilterDatabase")
$fd.RefersTo = "=Planilha1!`$A`$1:`$O`$1"
